$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.123.50"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "1.589.16"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("D4").Value = "0.9960"
$ws.Range("E4").Value = "  -0.55%  "

$ws.Range("D5").Value = "0.9984"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").Value = "300.76"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "0.3765"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "50.95"
$ws.Range("E8").Value = "  +4.86%  "

$ws.Range("D9").Value = "0.3593"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").Value = "1.236"
$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("D11").Value = "0.9992"
$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").Value = "0.08029"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").Value = "22.24"
$ws.Range("E13").Value = "  -2.51%  "

$ws.Range("D14").Value = "6.496"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "7.341"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "0.00001238"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").Value = "1.586.92"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "92.84"
$ws.Range("E18").Value = "  +1.84%  "

$ws.Range("D19").Value = "0.06756"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "17.91"
$ws.Range("E20").Value = "  -1.84%  "

$ws.Range("D21").Value = "0.9989"
$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").Value = "6.398"
$ws.Range("E22").Value = "  -2.04%  "

$ws.Range("D23").Value = "23.042.40"
$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("E24").Value = "  -2.40%  "

$ws.Range("D25").Value = "2.380"
$ws.Range("E25").Value = "  +1.38%  "

$ws.Range("D26").Value = "2.861"
$ws.Range("E26").Value = "  +2.37%  "

$ws.Range("D27").Value = "20.79"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("D28").Value = "148.18"
$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("D29").Value = "5.188"
$ws.Range("E29").Value = "  -1.15%  "

$ws.Range("D30").Value = "132.36"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("D31").Value = "2.358"
$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("D32").Value = "6.615"
$ws.Range("E32").Value = "  -2.53%  "

$ws.Range("D33").Value = "1.762.91"
$ws.Range("E33").Value = "  -0.71%  "

$ws.Range("D34").Value = "0.9511"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").Value = "0.07436"
$ws.Range("E35").Value = "  -2.77%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.02667"
$ws.Range("E36").Value = "  -3.30%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "10.02"
$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("D38").Value = "0.08782"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("D39").Value = "0.2491"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").Value = "6.078"
$ws.Range("E40").Value = "  -1.84%  "

$ws.Range("D41").Value = "1.351"
$ws.Range("E41").Value = "  -2.32%  "

$ws.Range("D42").Value = "0.6968"
$ws.Range("E42").Value = "  -2.45%  "

$ws.Range("D43").Value = "12.16"
$ws.Range("E43").Value = "  -4.48%  "

$ws.Range("D44").Value = "14.93"
$ws.Range("E44").Value = "  -4.11%  "

$ws.Range("D45").Value = "0.6449"
$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("D46").Value = "0.9979"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("D49").Value = "131.19"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("D50").Value = "0.07878"
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("D51").Value = "1.202"
$ws.Range("E51").Value = "  +3.04%  "
